$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 45 - new date (2023-01-06), reuse the date style from A39 (s="1")
$ws.Range("A39").Copy()
$ws.Range("A45").PasteSpecial(-4122)
$ws.Range("A45").Value = 44932

# Row 46 - time + note, reuse the time style from A40 (s="2")
$ws.Range("A40").Copy()
$ws.Range("A46").PasteSpecial(-4122)
$ws.Range("A46").Value = 0.37847222222222227
$ws.Range("B46").Value = "Switch"

# Row 47 - time + note
$ws.Range("A40").Copy()
$ws.Range("A47").PasteSpecial(-4122)
$ws.Range("A47").Value = 0.43611111111111112
$ws.Range("B47").Value = "code產生器"

# Row 48 - time + note
$ws.Range("A40").Copy()
$ws.Range("A48").PasteSpecial(-4122)
$ws.Range("A48").Value = 0.49722222222222223
$ws.Range("B48").Value = "const v.s. readonly"

# Row 49 - column D note (written before D48 so shared-string order matches)
$ws.Range("D49").Value = "const會在編譯時直接替換成值"

# Row 48 - column D note
$ws.Range("D48").Value = "Readonly變數只能在建構子給值"

# Update selection to match the saved view state
$ws.Range("H41").Select()
